$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 308-309; this shifts the former rows
# 308..328 down to 310..330 and extends the sheet dimension to R330.
$ws.Range("A308:R309").EntireRow.Insert()

# --- New row 308 ---
$ws.Cells.Item(308, 1).Value = 8
$ws.Cells.Item(308, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(308, 3).Value = "Coquimbo"
$ws.Cells.Item(308, 4).Value = 44714
$ws.Cells.Item(308, 5).Value = 4
$ws.Cells.Item(308, 6).Value = 100112032
$ws.Cells.Item(308, 7).Value = "Zapallo italiano"
$ws.Cells.Item(308, 8).Value = "Bola 8"
$ws.Cells.Item(308, 9).Value = "Primera"
$ws.Cells.Item(308, 10).Value = 400
$ws.Cells.Item(308, 11).Value = 10000
$ws.Cells.Item(308, 12).Value = 11000
$ws.Cells.Item(308, 13).Value = 10500
$ws.Cells.Item(308, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(308, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(308, 16).Value = 210
$ws.Cells.Item(308, 17).Value = 50
$ws.Cells.Item(308, 18).Value = "Hortaliza"

# --- New row 309 ---
$ws.Cells.Item(309, 1).Value = 8
$ws.Cells.Item(309, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 44714
$ws.Cells.Item(309, 5).Value = 4
$ws.Cells.Item(309, 6).Value = 100112032
$ws.Cells.Item(309, 7).Value = "Zapallo italiano"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 540
$ws.Cells.Item(309, 11).Value = 13000
$ws.Cells.Item(309, 12).Value = 14000
$ws.Cells.Item(309, 13).Value = 13500
$ws.Cells.Item(309, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(309, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(309, 16).Value = 225
$ws.Cells.Item(309, 17).Value = 60
$ws.Cells.Item(309, 18).Value = "Hortaliza"

# Give the two new date cells the same number format as the rest of
# column D (style index 2 / numFmtId 165 date-time format).
$ws.Range("D308:D309").NumberFormat = $ws.Range("D310").NumberFormat
